$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.275.87"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "3.321.58"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "248.22"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("D6").Value = "650.66"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -6.54%  "
$ws.Range("D8").Value = "0.418"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -7.18%  "
$ws.Range("D11").Value = "3.318.04"
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("E12").Value = "  -3.36%  "
$ws.Range("D13").Value = "40.03"
$ws.Range("E13").Value = "  -4.54%  "
$ws.Range("D14").Value = "95.995.16"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "6.05"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("E16").Value = "  -3.90%  "
$ws.Range("D17").Value = "3.939.38"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "3.338.05"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("E20").Value = "  +3.51%  "
$ws.Range("D21").Value = "16.99"
$ws.Range("E21").Value = "  -2.67%  "
$ws.Range("D22").Value = "502.24"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").Value = "  -4.13%  "
$ws.Range("D24").Value = "3.35"
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("E25").Value = "  -4.32%  "
$ws.Range("E26").Value = "  +6.62%  "
$ws.Range("D27").Value = "95.66"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("Z1").Formula = "=""12.00"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E28").Value = "  -5.74%  "
$ws.Range("D29").Value = "0.143"
$ws.Range("E29").Value = "  -8.12%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  -3.77%  "
$ws.Range("D32").Value = "0.187"
$ws.Range("E32").Value = "  -6.33%  "
$ws.Range("D33").Value = "2.46"
$ws.Range("E33").Value = "  +8.37%  "
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").Value = "0.543"
$ws.Range("E35").Value = "  -5.47%  "
$ws.Range("D36").Value = "27.87"
$ws.Range("E36").Value = "  -6.67%  "
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").Value = "503.15"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").Value = "24.34"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").Value = "0.0428"
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("D44").Value = "0.825"
$ws.Range("E44").Value = "  -3.66%  "
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("E46").Value = "  +5.46%  "
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("D49").Value = "52.96"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("E50").Value = "  -5.20%  "
$ws.Range("D51").Value = "161.84"
$ws.Range("E51").Value = "  -0.05%  "
